$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "Leer archivo de texto"
$ws.Range("E4").Value = "Permite la lectura de un archivo de texto con la información de los asistentes (espectadores y posibles participantes)"
$ws.Range("F4").Value = "Ruta del archivo"
$ws.Range("G4").Value = "Se lee el archivo de texto con la información correctamente "

$ws.Range("D5").Value = "Cargar información"
$ws.Range("E5").Value = "Permite cargar la información leída del archivo de texto en un árbol binario"

$ws.Rows.Item(4).RowHeight = 43.2
$ws.Rows.Item(5).RowHeight = 28.8

$ws.Range("F5").Select()
